$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new column before column A. Excel's native "insert column"
# shifts every existing column (A..AC -> B..AD), their values, styles,
# merged cells, and the used-range dimension along with it.
$ws.Columns("A").Insert()

# Populate the new column A with the "Match ID" field.
# Row 3 is the (visible) header row -> label the new column.
$ws.Range("A3").Value = "Match ID"
$ws.Range("A3").Font.Bold = $true

# Rows 4-19 are the visible per-player data rows -> constant match id (33),
# styled with the bold header font (no border/fill) like the rest of col A.
for ($r = 4; $r -le 19; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = 33
    $cell.Font.Bold = $true
}

# Row 20 is the hidden totals/summary row -> same match id, default style.
$ws.Cells.Item(20, 1).Value = 33
# Re-fit the row height so writing into the hidden row doesn't stamp an
# explicit custom row height onto it.
$ws.Rows(20).AutoFit()

# Match the saved selection state (A3:A19 selected, active cell A3).
$ws.Range("A3:A19").Select() | Out-Null
